$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.40130033333333
$ws.Range("H2").Value = 109.203901
$ws.Range("I2").Value = 0.1897437225523226
$ws.Range("J2").Value = 0.1897437225523226
$ws.Range("M2").Value = 7.021661333333333
$ws.Range("N2").Value = 21.064984
$ws.Range("O2").Value = 0.5944820341692109
$ws.Range("P2").Value = 0.5944820341692109
$ws.Range("Q2").Value = 255.5976030336204
$ws.Range("R2").Value = 2300.378427302584
$ws.Range("S2").Value = 0.1127992341537431
$ws.Range("T2").Value = 0.1127992341537431
$ws.Range("G3").Value = 36.40130033333333
$ws.Range("H3").Value = 109.203901
$ws.Range("I3").Value = 0.1897437225523226
$ws.Range("J3").Value = 0.1897437225523226
$ws.Range("M3").Value = 7.021661333333333
$ws.Range("N3").Value = 21.064984
$ws.Range("O3").Value = 0.09094063723386185
$ws.Range("P3").Value = 0.09094063723386187
$ws.Range("Q3").Value = 39.09993500107855
$ws.Range("R3").Value = 351.899415009707
$ws.Range("S3").Value = 0.0172554150400333
$ws.Range("T3").Value = 0.01725541504003331
$ws.Range("G4").Value = 36.40130033333333
$ws.Range("H4").Value = 109.203901
$ws.Range("I4").Value = 0.1897437225523226
$ws.Range("J4").Value = 0.1897437225523226
$ws.Range("M4").Value = 3.715596666666666
$ws.Range("N4").Value = 11.14679
$ws.Range("O4").Value = 0.3145773285969274
$ws.Range("P4").Value = 0.3145773285969274
$ws.Range("Q4").Value = 135.2525501808655
$ws.Range("R4").Value = 1217.27295162779
$ws.Range("S4").Value = 0.0596890733585462
$ws.Range("T4").Value = 0.05968907335854621
$ws.Range("I5").Value = 0.6107553255746098
$ws.Range("J5").Value = 0.6107553255746098
$ws.Range("M5").Value = 7.021661333333333
$ws.Range("N5").Value = 21.064984
$ws.Range("O5").Value = 0.5944820341692109
$ws.Range("P5").Value = 0.5944820341692109
$ws.Range("Q5").Value = 822.7286529273264
$ws.Range("R5").Value = 7404.557876345936
$ws.Range("S5").Value = 0.3630830683272727
$ws.Range("T5").Value = 0.3630830683272727
$ws.Range("I6").Value = 0.6107553255746098
$ws.Range("J6").Value = 0.6107553255746098
$ws.Range("O6").Value = 0.09094063723386185
$ws.Range("P6").Value = 0.09094063723386187
$ws.Range("S6").Value = 0.05554247850172978
$ws.Range("T6").Value = 0.05554247850172978
$ws.Range("I7").Value = 0.6107553255746098
$ws.Range("J7").Value = 0.6107553255746098
$ws.Range("M7").Value = 3.715596666666666
$ws.Range("N7").Value = 11.14679
$ws.Range("O7").Value = 0.3145773285969274
$ws.Range("P7").Value = 0.3145773285969274
$ws.Range("Q7").Value = 435.3567760205178
$ws.Range("R7").Value = 3918.21098418466
$ws.Range("S7").Value = 0.1921297787456074
$ws.Range("T7").Value = 0.1921297787456074
$ws.Range("G8").Value = 38.27317166666666
$ws.Range("H8").Value = 114.819515
$ws.Range("I8").Value = 0.1995009518730676
$ws.Range("J8").Value = 0.1995009518730676
$ws.Range("M8").Value = 7.021661333333333
$ws.Range("N8").Value = 21.064984
$ws.Range("O8").Value = 0.5944820341692109
$ws.Range("P8").Value = 0.5944820341692109
$ws.Range("Q8").Value = 268.7412495958622
$ws.Range("R8").Value = 2418.67124636276
$ws.Range("S8").Value = 0.1185997316881951
$ws.Range("T8").Value = 0.1185997316881951
$ws.Range("G9").Value = 38.27317166666666
$ws.Range("H9").Value = 114.819515
$ws.Range("I9").Value = 0.1995009518730676
$ws.Range("J9").Value = 0.1995009518730676
$ws.Range("M9").Value = 7.021661333333333
$ws.Range("N9").Value = 21.064984
$ws.Range("O9").Value = 0.09094063723386185
$ws.Range("P9").Value = 0.09094063723386187
$ws.Range("Q9").Value = 41.11057876362277
$ws.Range("R9").Value = 369.995208872605
$ws.Range("S9").Value = 0.01814274369209877
$ws.Range("T9").Value = 0.01814274369209878
$ws.Range("G10").Value = 38.27317166666666
$ws.Range("H10").Value = 114.819515
$ws.Range("I10").Value = 0.1995009518730676
$ws.Range("J10").Value = 0.1995009518730676
$ws.Range("M10").Value = 3.715596666666666
$ws.Range("N10").Value = 11.14679
$ws.Range("O10").Value = 0.3145773285969274
$ws.Range("P10").Value = 0.3145773285969274
$ws.Range("Q10").Value = 142.2076690674278
$ws.Range("R10").Value = 1279.86902160685
$ws.Range("S10").Value = 0.06275847649277379
$ws.Range("T10").Value = 0.06275847649277379